$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(6,7,9,15,22,24,27,29,31,34,35,36,37,39,41,45,52,54,55,58,64,65,67,72,75,76,79,81,82,85,91,92,98,100,102,105,107,109,110,111,115,119,120,123,129,136,138,139,140,144,150,157,159,161,162,167,170,177,178,179,181,186,188,193)

foreach ($r in $rows) {
    $ws.Range("G$r`:N$r").Value = "N/A"
}
